# Adds the I0 and IF columns (I and J) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy formatting from an existing header cell (H1) so the
# new headers match the rest of the header row (bold, bordered, centered).
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$data = @(
    @(7,8), @(5,6), @(8,8), @(6,6), @(5,7), @(6,7), @(6,7), @(10,10), @(9,9), @(7,7),
    @(7,8), @(7,8), @(7,7), @(7,8), @(6,7), @(9,9), @(6,6), @(9,9), @(9,9), @(7,8),
    @(6,6), @(5,7), @(6,8), @(7,7), @(7,7), @(10,10), @(8,8), @(7,7), @(4,5), @(6,6),
    @(6,7), @(5,6), @(7,7), @(6,6), @(6,7), @(6,6), @(5,7), @(8,9), @(6,7), @(8,8),
    @(8,9), @(1,4), @(1,4), @(1,3), @(7,9), @(3,4)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
